# Slide 1, "TextBox 3" shape (3rd shape on the slide) lists team members,
# one per paragraph:
#   1. Alex
#   2. Conner
#   3. Theodore Sosnowski
#   4. Chris Wilson
#   5. (blank)
#
# This edit:
#   - Para 1: "Alex" -> "Alex " + "Wimer" (adds a last name as a new run)
#   - Para 2: "Conner" -> "Connor Mahaffey" (fixes the spelling / adds surname)

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)
$tr  = $shp.TextFrame.TextRange

# --- Paragraph 1: "Alex" -> "Alex " + "Wimer" --------------------------
$para1 = $tr.Paragraphs(1, 1)
$run1  = $para1.Runs(1)
$run1.Text = "Alex "
[void]$run1.InsertAfter("Wimer")

# --- Paragraph 2: "Conner" -> "Connor Mahaffey" -------------------------
$para2 = $tr.Paragraphs(2, 1)
$run2  = $para2.Runs(1)
$run2.Text = "Connor Mahaffey"
